$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force the Price column to Text so numeric-looking values
# (e.g. "206.46") are written as strings, matching the source inline-string cells.
$dRange = $ws.Range("D2:D50")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = '27.680.67'
$ws.Range("D3").Value = '1.584.92'
$ws.Range("D5").Value = '206.46'
$ws.Range("D8").Value = '22.25'
$ws.Range("D12").Value = '1.809.94'
$ws.Range("D13").Value = '1.568.67'
$ws.Range("D16").Value = '27.651.42'
$ws.Range("D17").Value = '63.28'
$ws.Range("D18").Value = '220.60'
$ws.Range("D20").Value = '7.31'
$ws.Range("D23").Value = '9.48'
$ws.Range("D25").Value = '153.95'
$ws.Range("D28").Value = '15.12'
$ws.Range("D31").Value = '0.0464'
$ws.Range("D33").Value = '1.386.88'
$ws.Range("D39").Value = '0.540'
$ws.Range("D42").Value = '0.977'
$ws.Range("D45").Value = '63.55'
$ws.Range("D47").Value = '1.721.66'
$ws.Range("D48").Value = '87.97'
$ws.Range("D50").Value = '0.0974'

# Restore the default (unstyled) formatting now that the text values are set
$dRange.Style = "Normal"

# Volume(1h) column updates
$ws.Range("E2").Value = '  -0.87%  '
$ws.Range("E3").Value = '  -3.13%  '
$ws.Range("E4").Value = '  +0.31%  '
$ws.Range("E5").Value = '  -2.54%  '
$ws.Range("E6").Value = '  -2.88%  '
$ws.Range("E7").Value = '  +0.37%  '
$ws.Range("E8").Value = '  -4.85%  '
$ws.Range("E9").Value = '  -1.32%  '
$ws.Range("E10").Value = '  -3.16%  '
$ws.Range("E11").Value = '  -1.84%  '
$ws.Range("E12").Value = '  -3.11%  '
$ws.Range("E13").Value = '  -4.12%  '
$ws.Range("E14").Value = '  -4.11%  '
$ws.Range("E15").Value = '  -5.86%  '
$ws.Range("E16").Value = '  -0.96%  '
$ws.Range("E17").Value = '  -3.11%  '
$ws.Range("E18").Value = '  -3.74%  '
$ws.Range("E19").Value = '  -3.78%  '
$ws.Range("E20").Value = '  -5.51%  '
$ws.Range("E21").Value = '  +0.33%  '
$ws.Range("E22").Value = '  -5.34%  '
$ws.Range("E23").Value = '  -6.70%  '
$ws.Range("E24").Value = '  -5.85%  '
$ws.Range("E25").Value = '  -1.31%  '
$ws.Range("E26").Value = '  -2.77%  '
$ws.Range("E27").Value = '  +0.25%  '
$ws.Range("E28").Value = '  -2.86%  '
$ws.Range("E29").Value = '  -4.17%  '
$ws.Range("E30").Value = '  -2.71%  '
$ws.Range("E31").Value = '  -3.53%  '
$ws.Range("E32").Value = '  -5.83%  '
$ws.Range("E33").Value = '  -0.96%  '
$ws.Range("E34").Value = '  -5.38%  '
$ws.Range("E35").Value = '  -5.27%  '
$ws.Range("E36").Value = '  -5.21%  '
$ws.Range("E37").Value = '  -0.94%  '
$ws.Range("E38").Value = '  -3.06%  '
$ws.Range("E39").Value = '  -3.49%  '
$ws.Range("E40").Value = '  -3.86%  '
$ws.Range("E41").Value = '  +0.29%  '
$ws.Range("E42").Value = '  -3.08%  '
$ws.Range("E43").Value = '  -4.22%  '
$ws.Range("E44").Value = '  +1.63%  '
$ws.Range("E45").Value = '  -3.85%  '
$ws.Range("E46").Value = '  -4.06%  '
$ws.Range("E47").Value = '  -3.06%  '
$ws.Range("E48").Value = '  -0.87%  '
$ws.Range("E49").Value = '  -2.50%  '
$ws.Range("E50").Value = '  -5.04%  '
